$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 211-213 (shifts existing rows 211-289 down to 214-292)
$ws.Range("A211:A213").EntireRow.Insert()

# Row 211: new Venus entry
$ws.Cells.Item(211,1).Value = 11
$ws.Cells.Item(211,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(211,3).Value = "Bíobío"
$ws.Cells.Item(211,4).Value = 44609
$ws.Cells.Item(211,5).Value = 8
$ws.Cells.Item(211,6).Value = "Fruta"
$ws.Cells.Item(211,7).Value = 100103
$ws.Cells.Item(211,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(211,9).Value = 100103006
$ws.Cells.Item(211,10).Value = "Nectarín"
$ws.Cells.Item(211,11).Value = "Venus"
$ws.Cells.Item(211,12).Value = "Especial"
$ws.Cells.Item(211,13).Value = 50
$ws.Cells.Item(211,14).Value = 13000
$ws.Cells.Item(211,15).Value = 13000
$ws.Cells.Item(211,16).Value = 13000
$ws.Cells.Item(211,17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(211,18).Value = "Región de O'Higgins"
$ws.Cells.Item(211,19).Value = 812
$ws.Cells.Item(211,20).Value = 16

# Row 212: new Venus entry
$ws.Cells.Item(212,1).Value = 11
$ws.Cells.Item(212,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(212,3).Value = "Bíobío"
$ws.Cells.Item(212,4).Value = 44609
$ws.Cells.Item(212,5).Value = 8
$ws.Cells.Item(212,6).Value = "Fruta"
$ws.Cells.Item(212,7).Value = 100103
$ws.Cells.Item(212,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(212,9).Value = 100103006
$ws.Cells.Item(212,10).Value = "Nectarín"
$ws.Cells.Item(212,11).Value = "Venus"
$ws.Cells.Item(212,12).Value = "Primera"
$ws.Cells.Item(212,13).Value = 100
$ws.Cells.Item(212,14).Value = 11000
$ws.Cells.Item(212,15).Value = 11000
$ws.Cells.Item(212,16).Value = 11000
$ws.Cells.Item(212,17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(212,18).Value = "Región de O'Higgins"
$ws.Cells.Item(212,19).Value = 688
$ws.Cells.Item(212,20).Value = 16

# Row 213: new Venus entry
$ws.Cells.Item(213,1).Value = 11
$ws.Cells.Item(213,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(213,3).Value = "Bíobío"
$ws.Cells.Item(213,4).Value = 44609
$ws.Cells.Item(213,5).Value = 8
$ws.Cells.Item(213,6).Value = "Fruta"
$ws.Cells.Item(213,7).Value = 100103
$ws.Cells.Item(213,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(213,9).Value = 100103006
$ws.Cells.Item(213,10).Value = "Nectarín"
$ws.Cells.Item(213,11).Value = "Venus"
$ws.Cells.Item(213,12).Value = "Segunda"
$ws.Cells.Item(213,13).Value = 100
$ws.Cells.Item(213,14).Value = 9000
$ws.Cells.Item(213,15).Value = 9000
$ws.Cells.Item(213,16).Value = 9000
$ws.Cells.Item(213,17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(213,18).Value = "Región de O'Higgins"
$ws.Cells.Item(213,19).Value = 562
$ws.Cells.Item(213,20).Value = 16
